$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 237.17647
$ws.Range("I33").Value = 226.73334
$ws.Range("K33").Value = 226.73334
$ws.Range("M33").Value = 2.266660000000002
$ws.Range("H82").Value = 9884.286
$ws.Range("I82").Value = 2621
$ws.Range("J82").Value = 15331.75
$ws.Range("K82").Value = 7863
$ws.Range("L82").Value = 45995.25
$ws.Range("M82").Value = -7457
$ws.Range("N82").Value = -46807.25
$ws.Range("H85").Value = 9884.286
$ws.Range("I85").Value = 2621
$ws.Range("J85").Value = 15331.75
$ws.Range("K85").Value = 7863
$ws.Range("L85").Value = 45995.25
$ws.Range("M85").Value = -6459
$ws.Range("N85").Value = -48803.25
$ws.Range("H101").Value = 1536.875
$ws.Range("I101").Value = 178.2
$ws.Range("J101").Value = 3801.3333
$ws.Range("K101").Value = 534.5999999999999
$ws.Range("L101").Value = 11403.9999
$ws.Range("M101").Value = 1087.4
$ws.Range("N101").Value = -14647.9999
$ws.Range("H137").Value = 26321008
$ws.Range("I137").Value = 55559410
$ws.Range("K137").Value = 166678230
$ws.Range("M137").Value = -166675680
$ws.Range("H138").Value = 5825.3335
$ws.Range("I138").Value = 3634.04
$ws.Range("J138").Value = 9047.823
$ws.Range("K138").Value = 10902.12
$ws.Range("L138").Value = 27143.469
$ws.Range("M138").Value = -5762.119999999999
$ws.Range("N138").Value = -37423.469
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4728.5957
$ws.Range("I32").Value = 4427.5347
$ws.Range("K32").Value = 4427.5347
$ws.Range("M32").Value = -4140.5347
$ws.Range("H43").Value = 24256.334
$ws.Range("J43").Value = 19988.143
$ws.Range("L43").Value = 19988.143
$ws.Range("N43").Value = -20614.143
$ws.Range("H55").Value = 23399.6
$ws.Range("H74").Value = 2014.5312
$ws.Range("I74").Value = 1739.7037
$ws.Range("K74").Value = 1739.7037
$ws.Range("M74").Value = -865.7037
$ws.Range("H77").Value = 2014.5312
$ws.Range("I77").Value = 1739.7037
$ws.Range("K77").Value = 8698.5185
$ws.Range("M77").Value = -4330.5185
$ws.Range("H88").Value = 2136.375
$ws.Range("J88").Value = 2588.5
$ws.Range("L88").Value = 2588.5
$ws.Range("N88").Value = -3400.5
$ws.Range("H91").Value = 2136.375
$ws.Range("J91").Value = 2588.5
$ws.Range("L91").Value = 2588.5
$ws.Range("N91").Value = -5396.5
$ws.Range("H97").Value = 1283.3334
$ws.Range("I97").Value = 1165.8334
$ws.Range("J97").Value = 1596.6666
$ws.Range("K97").Value = 1165.8334
$ws.Range("L97").Value = 1596.6666
$ws.Range("M97").Value = -669.8334
$ws.Range("N97").Value = -2588.6666
$ws.Range("H102").Value = 3417.6
$ws.Range("I102").Value = 2685.1428
$ws.Range("K102").Value = 2685.1428
$ws.Range("M102").Value = -1063.1428
$ws.Range("H110").Value = 6432.727
$ws.Range("I110").Value = 7545.125
$ws.Range("K110").Value = 7545.125
$ws.Range("M110").Value = -5500.125
$ws.Range("H132").Value = 2384820
$ws.Range("I132").Value = 3810.975
$ws.Range("J132").Value = 50005000
$ws.Range("K132").Value = 11432.925
$ws.Range("L132").Value = 150015000
$ws.Range("M132").Value = -8902.924999999999
$ws.Range("N132").Value = -150020060
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4764670
$ws.Range("I134").Value = 2361.2856
$ws.Range("K134").Value = 7083.8568
$ws.Range("M134").Value = -4548.8568
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1061.2222
$ws.Range("I22").Value = 793.875
$ws.Range("K22").Value = 793.875
$ws.Range("M22").Value = -443.875
$ws.Range("H59").Value = 118799.2
$ws.Range("J59").Value = 142999.5
$ws.Range("L59").Value = 142999.5
$ws.Range("N59").Value = -145289.5
$ws.Range("H102").Value = 99120
$ws.Range("J102").Value = 99120
$ws.Range("L102").Value = 99120
$ws.Range("N102").Value = -103988
$ws.Range("H124").Value = 79599.336
$ws.Range("J124").Value = 79599.336
$ws.Range("L124").Value = 79599.336
$ws.Range("N124").Value = -84509.336
$ws.Range("H134").Value = 4788.5
$ws.Range("I134").Value = 4343.6665
$ws.Range("J134").Value = 8124.75
$ws.Range("K134").Value = 13030.9995
$ws.Range("L134").Value = 24374.25
$ws.Range("M134").Value = -10495.9995
$ws.Range("N134").Value = -29444.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1668
$ws.Range("I5").Value = 1576.5
$ws.Range("J5").Value = 1741.2
$ws.Range("K5").Value = 4729.5
$ws.Range("L5").Value = 5223.6
$ws.Range("M5").Value = -4617.5
$ws.Range("N5").Value = -5447.6
$ws.Range("H68").Value = 2308.9333
$ws.Range("I68").Value = 1700
$ws.Range("J68").Value = 2352.4285
$ws.Range("K68").Value = 5100
$ws.Range("L68").Value = 7057.2855
$ws.Range("M68").Value = -4289
$ws.Range("N68").Value = -8679.2855
$ws.Range("H71").Value = 2308.9333
$ws.Range("I71").Value = 1700
$ws.Range("J71").Value = 2352.4285
$ws.Range("K71").Value = 15300
$ws.Range("L71").Value = 21171.8565
$ws.Range("M71").Value = -11244
$ws.Range("N71").Value = -29283.8565
$ws.Range("H97").Value = 1739.625
$ws.Range("I97").Value = 1948.6666
$ws.Range("J97").Value = 1614.2
$ws.Range("K97").Value = 5845.9998
$ws.Range("L97").Value = 4842.6
$ws.Range("M97").Value = -5349.9998
$ws.Range("N97").Value = -5834.6
$ws.Range("H135").Value = 1668
$ws.Range("I135").Value = 1576.5
$ws.Range("J135").Value = 1741.2
$ws.Range("K135").Value = 14188.5
$ws.Range("L135").Value = 15670.8
$ws.Range("M135").Value = -11653.5
$ws.Range("N135").Value = -20740.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 60001576
$ws.Range("I80").Value = 1716.875
$ws.Range("J80").Value = 300001000
$ws.Range("K80").Value = 1716.875
$ws.Range("L80").Value = 300001000
$ws.Range("M80").Value = -718.875
$ws.Range("N80").Value = -300002996
$ws.Range("H83").Value = 60001576
$ws.Range("I83").Value = 1716.875
$ws.Range("J83").Value = 300001000
$ws.Range("K83").Value = 8584.375
$ws.Range("L83").Value = 1500005000
$ws.Range("M83").Value = -3592.375
$ws.Range("N83").Value = -1500014984
$ws.Range("H102").Value = 2561.85
$ws.Range("I102").Value = 2356.8462
$ws.Range("K102").Value = 2356.8462
$ws.Range("M102").Value = -734.8462
$ws.Range("H123").Value = 96306
$ws.Range("J123").Value = 96306
$ws.Range("L123").Value = 96306
$ws.Range("N123").Value = -101206
$ws.Range("H132").Value = 2002381.5
$ws.Range("I132").Value = 2410.175
$ws.Range("K132").Value = 7230.525000000001
$ws.Range("M132").Value = -4700.525000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2232.5938
$ws.Range("I40").Value = 2175.6128
$ws.Range("K40").Value = 2175.6128
$ws.Range("M40").Value = -2039.6128
$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("K46").Value = 500
$ws.Range("M46").Value = -312
$ws.Range("H55").Value = 1048.591
$ws.Range("I55").Value = 519.8
$ws.Range("J55").Value = 1489.25
$ws.Range("K55").Value = 519.8
$ws.Range("L55").Value = 1489.25
$ws.Range("M55").Value = -346.8
$ws.Range("N55").Value = -1835.25
$ws.Range("H68").Value = 2317324.2
$ws.Range("J68").Value = 2686.125
$ws.Range("L68").Value = 2686.125
$ws.Range("N68").Value = -4184.125
$ws.Range("H71").Value = 2317324.2
$ws.Range("J71").Value = 2686.125
$ws.Range("L71").Value = 13430.625
$ws.Range("N71").Value = -20918.625
$ws.Range("H132").Value = 4526.067
$ws.Range("I132").Value = 3049.1428
$ws.Range("K132").Value = 9147.428400000001
$ws.Range("M132").Value = -6617.428400000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3998.8333
$ws.Range("I122").Value = 3998.25
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 11994.75
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -9544.75
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 279503.12
$ws.Range("I132").Value = 1562.3226
$ws.Range("J132").Value = 2002736.2
$ws.Range("K132").Value = 4686.9678
$ws.Range("L132").Value = 6008208.6
$ws.Range("M132").Value = -2156.9678
$ws.Range("N132").Value = -6013268.6
$ws.Range("H136").Value = 226845.83
$ws.Range("I136").Value = 14470.738
$ws.Range("J136").Value = 2010796.6
$ws.Range("K136").Value = 43412.214
$ws.Range("L136").Value = 6032389.800000001
$ws.Range("M136").Value = -40862.214
$ws.Range("N136").Value = -6037489.800000001
